# Rename the "质控组" (QC group) to "北京组" (Beijing group) across the
# dashboard workbook: Sheet1 column A (rows 2-5) and Sheet2 cell A2.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("A2").Value = "北京组"
$ws1.Range("A3").Value = "北京组"
$ws1.Range("A4").Value = "北京组"
$ws1.Range("A5").Value = "北京组"

$ws2.Range("A2").Value = "北京组"

# Reflect the editor's final focus: Sheet2 was reviewed (A2 selected) and
# then Sheet1 was left as the active sheet with A5 selected.
$ws2.Activate()
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("A5").Select()
